$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Staging cell used purely to coerce the new Order ID values into being
# stored as text (shared string) while leaving the destination cell's
# own number format (General) untouched - mirrors how the original
# numeric-looking OrderId values were stored as text in this column.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"

$updates = @{
    "R2" = "51532508"
    "R3" = "51532255"
    "R4" = "51532510"
    "R5" = "51532511"
}

foreach ($addr in $updates.Keys) {
    $helper.Value = $updates[$addr]
    $helper.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

$helper.Clear()
$excel.CutCopyMode = 0
